# [ENote] Small change in api map
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Simplify the search endpoint paths (drop the {title}/{label} placeholder suffix)
$ws.Range("A2").Value = "notes/search_title"
$ws.Range("A3").Value = "notes/search_label"

# Both search endpoints are actually POST requests, not GET
$ws.Range("B2").Value = "POST"
$ws.Range("B3").Value = "POST"

# Move the active selection from A5 to D4
$ws.Range("D4").Select()
